# Weekly price-data refresh for "Hortaliza, Agrícola del Norte S.A. de Arica - Pepino dulce"
# - Revises the figures on existing rows 7-15 (new market week shifted the series down)
# - Appends three additional observation rows (16-18) that were pushed out of the old range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ "A" = 1; "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6; "G" = 7; "H" = 8; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14; "O" = 15; "P" = 16; "Q" = 17; "R" = 18 }

# Cells that change value on the already-existing rows (7-15)
$updates = @(
    @{ Row = 7; Cells = @{ "D" = 44526; "H" = "Cultivar XV región"; "I" = "Primera"; "J" = 100; "K" = 5000; "L" = 5500; "M" = 5250; "N" = "`$/caja 10 kilos"; "O" = "Región de Arica y Parinacota"; "P" = 525; "Q" = 10 } },
    @{ Row = 8; Cells = @{ "D" = 44526; "H" = "Cultivar XV región"; "K" = 4000; "L" = 4500; "M" = 4250; "N" = "`$/caja 10 kilos"; "O" = "Región de Arica y Parinacota"; "P" = 425; "Q" = 10 } },
    @{ Row = 9; Cells = @{ "D" = 44526; "H" = "Cultivar XV región"; "K" = 3000; "L" = 3500; "M" = 3250; "N" = "`$/caja 10 kilos"; "O" = "Región de Arica y Parinacota"; "P" = 325; "Q" = 10 } },
    @{ Row = 10; Cells = @{ "D" = 44405; "I" = "Segunda"; "J" = 140 } },
    @{ Row = 11; Cells = @{ "D" = 44435; "K" = 17000; "L" = 18000; "M" = 17500; "P" = 972 } },
    @{ Row = 12; Cells = @{ "D" = 44435; "H" = "Cultivar IV Región"; "I" = "Tercera"; "J" = 120; "K" = 14000; "L" = 15000; "M" = 14500; "N" = "`$/bandeja 18 kilos"; "O" = "Provincia de Limarí"; "P" = 806; "Q" = 18 } },
    @{ Row = 13; Cells = @{ "D" = 44398; "J" = 100; "K" = 17000; "L" = 18000; "M" = 17500; "P" = 972 } },
    @{ Row = 14; Cells = @{ "D" = 44398; "I" = "Segunda"; "J" = 100; "K" = 15000; "L" = 16000; "M" = 15500; "P" = 861 } },
    @{ Row = 15; Cells = @{ "D" = 44211; "H" = "Cultivar XV región"; "J" = 140; "K" = 4500; "L" = 5000; "M" = 4750; "N" = "`$/caja 10 kilos"; "O" = "Región de Arica y Parinacota"; "P" = 475; "Q" = 10 } }
)

# Brand-new rows appended at the bottom of the table (16-18)
$newRows = @(
    @{ Row = 16; Cells = @{ "A" = 1; "B" = "Agrícola del Norte S.A. de Arica"; "C" = "Arica y Parinacota"; "D" = 44454; "E" = 15; "F" = 100112043; "G" = "Pepino dulce"; "H" = "Cultivar IV Región"; "I" = "Primera"; "J" = 160; "K" = 19000; "L" = 20000; "M" = 19500; "N" = "`$/bandeja 18 kilos"; "O" = "Provincia de Limarí"; "P" = 1083; "Q" = 18; "R" = "Hortaliza" } },
    @{ Row = 17; Cells = @{ "A" = 1; "B" = "Agrícola del Norte S.A. de Arica"; "C" = "Arica y Parinacota"; "D" = 44363; "E" = 15; "F" = 100112043; "G" = "Pepino dulce"; "H" = "Cultivar IV Región"; "I" = "Primera"; "J" = 140; "K" = 14000; "L" = 15000; "M" = 14500; "N" = "`$/bandeja 18 kilos"; "O" = "Provincia de Limarí"; "P" = 806; "Q" = 18; "R" = "Hortaliza" } },
    @{ Row = 18; Cells = @{ "A" = 1; "B" = "Agrícola del Norte S.A. de Arica"; "C" = "Arica y Parinacota"; "D" = 44391; "E" = 15; "F" = 100112043; "G" = "Pepino dulce"; "H" = "Cultivar IV Región"; "I" = "Segunda"; "J" = 100; "K" = 15000; "L" = 16000; "M" = 15500; "N" = "`$/bandeja 18 kilos"; "O" = "Provincia de Limarí"; "P" = 861; "Q" = 18; "R" = "Hortaliza" } }
)

foreach ($u in $updates) {
    foreach ($col in $u.Cells.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Cells[$col]
    }
}

# Template cell that already carries the correct date number format (style "s=2")
$dateFormat = $ws.Range("D2").NumberFormat

foreach ($u in $newRows) {
    foreach ($col in $u.Cells.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Cells[$col]
    }
    $ws.Cells.Item($u.Row, $colIndex["D"]).NumberFormat = $dateFormat
}
